# Change predict income algorithm
# Updates the forecasted income-statement figures (columns D:H, rows 11-27)
# on the "Overview" sheet, replacing placeholder zeros / "-" text with the
# newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Overview")

$values = @{
    11 = @(2213364, 3917739, 7208698, 11351159, 21157343)
    12 = @(-1344912, -1959861, -2952358, -5199910, -8863829)
    13 = @(868452, 1957878, 4256340, 6151249, 12293514)
    14 = @(-101106, -323609, -888297, -732994, -1806898)
    15 = @(0, 0, 0, 0, 0)
    16 = @(5429, 88662, 159893, -27547, 379865)
    17 = @(772775, 1722931, 3527936, 5390708, 10866481)
    18 = @(-36617, -44400, -43009, -48568, -55460)
    19 = @(132976, 142288, 275300, 454231, 1250178)
    20 = @(869134, 1820819, 3760227, 5796371, 12061199)
    21 = @(-169425, -307561, -396742, -689705, -1766521)
    22 = @(699709, 1513258, 3363485, 5106666, 10294678)
    23 = @(0, 0, 0, 0, 0)
    24 = @(699709, 1513258, 3363485, 5106666, 10294678)
    25 = @(1076, 0, 5175, 7856, 15838)
    26 = @(650000, 650000, 650000, 650000, 650000)
    27 = @(1076, 2328, 5175, 7856, 15838)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $col = 4 + $i   # column D is 4
        $ws.Cells.Item($row, $col).Value = $rowValues[$i]
    }
}

$wb.Save()
